$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (contaminants): B2 43->42, C2 1->3, D2 0.34->0.33, E2 0.01->0.02
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "42"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "3"
$ws.Range("D2").Value = 0.33
$ws.Range("E2").Value = 0.02

# Row 4 (flow_base_flow): B4 46->89, C4 13->7, D4 0.37->0.7, E4 0.1->0.06
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "89"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "7"
$ws.Range("D4").Value = 0.7
$ws.Range("E4").Value = 0.06

# Row 5 (food_web): C5 48->49, E5 0.38->0.39
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "49"
$ws.Range("E5").Value = 0.39

# Row 6 (temperature_rearing): B6 69->31, C6 15->7, D6 0.55->0.24, E6 0.12->0.06
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "31"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "7"
$ws.Range("D6").Value = 0.24
$ws.Range("E6").Value = 0.06
